$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$ws.Cells.Item(2, 12).Value = 0.9399999999999999
$ws.Cells.Item(3, 12).Value = 0.96
$ws.Cells.Item(4, 12).Value = 0.88
$ws.Cells.Item(5, 12).Value = 0.89
$ws.Cells.Item(6, 12).Value = 0.87
$ws.Cells.Item(7, 12).Value = 0.88
$ws.Cells.Item(8, 12).Value = 1.06
$ws.Cells.Item(9, 12).Value = 0.98
$ws.Cells.Item(10, 12).Value = 0.85
$ws.Cells.Item(11, 12).Value = 0.8100000000000001
$ws.Cells.Item(12, 12).Value = 0.8100000000000001
$ws.Cells.Item(13, 12).Value = 0.8
$ws.Cells.Item(14, 12).Value = 0.88
$ws.Cells.Item(15, 12).Value = 0.91
$ws.Cells.Item(16, 12).Value = 0.86
$ws.Cells.Item(17, 12).Value = 0.86
